$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $r = $ws.Range($range)
    $r.NumberFormat = "@"
    $r.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue "D2" "52.150.06"
Set-TextValue "E2" "  +1.20%  "

# Row 3 - Ethereum
Set-TextValue "D3" "2.870.32"
Set-TextValue "E3" "  +3.30%  "

# Row 4 - TetherUSD
Set-TextValue "E4" "  -0.04%  "

# Row 5 - BNB
Set-TextValue "D5" "350.98"
Set-TextValue "E5" "  +0.24%  "

# Row 6 - Solana
Set-TextValue "D6" "112.32"
Set-TextValue "E6" "  +3.84%  "

# Row 7 - XRP
Set-TextValue "D7" "0.561"
Set-TextValue "E7" "  +2.08%  "

# Row 8 - USDC
Set-TextValue "D8" "0.999"
Set-TextValue "E8" "  -0.03%  "

# Row 9 - Cardano
Set-TextValue "D9" "0.625"
Set-TextValue "E9" "  +2.28%  "

# Row 10 - Avalanche
Set-TextValue "D10" "40.26"
Set-TextValue "E10" "  +2.82%  "

# Row 11 - TRON
Set-TextValue "E11" "  +0.33%  "

# Row 12 - Dogecoin
Set-TextValue "D12" "0.0857"
Set-TextValue "E12" "  +3.05%  "

# Row 13 - Chainlink
Set-TextValue "D13" "20.07"
Set-TextValue "E13" "  +1.24%  "

# Row 14 - Polkadot
Set-TextValue "D14" "7.85"
Set-TextValue "E14" "  +1.31%  "

# Row 15 - WrappedliquidstakedEther2.0
Set-TextValue "D15" "3.319.26"
Set-TextValue "E15" "  +3.03%  "

# Row 16 - Polygon
Set-TextValue "D16" "0.994"
Set-TextValue "E16" "  +7.80%  "

# Row 17 - WrappedEther
Set-TextValue "D17" "2.876.73"
Set-TextValue "E17" "  +3.32%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "52.107.56"
Set-TextValue "E18" "  +1.24%  "

# Row 19 - ImmutableX
Set-TextValue "D19" "3.37"
Set-TextValue "E19" "  +9.65%  "

# Row 20 - Uniswap
Set-TextValue "D20" "7.68"
Set-TextValue "E20" "  -1.49%  "

# Row 21 - InternetComputer(DFINITY)
Set-TextValue "D21" "13.57"
Set-TextValue "E21" "  +2.55%  "

# Row 22 - ShibaInu
Set-TextValue "D22" "0.0₃0979"
Set-TextValue "E22" "  +1.80%  "

# Row 23 - Litecoin
Set-TextValue "D23" "70.78"
Set-TextValue "E23" "  +1.15%  "

# Row 24 - BitcoinCash
Set-TextValue "D24" "269.41"
Set-TextValue "E24" "  +1.31%  "

# Row 25 - PancakeSwap
Set-TextValue "E25" "  +1.44%  "

# Row 26 - EthereumClassic
Set-TextValue "D26" "26.36"
Set-TextValue "E26" "  +2.07%  "

# Row 27 - Dai
Set-TextValue "D27" "1.00"
Set-TextValue "E27" "  +0.17%  "

# Row 28 - Kaspa
Set-TextValue "D28" "0.165"
Set-TextValue "E28" "  +0.80%  "

# Row 29 - Cosmos
Set-TextValue "E29" "  +3.57%  "

# Row 30 - InjectiveProtocol
Set-TextValue "D30" "38.88"
Set-TextValue "E30" "  +4.73%  "

# Row 31 - Toncoin
Set-TextValue "E31" "  +1.18%  "

# Row 32 - Filecoin
Set-TextValue "E32" "  +1.62%  "

# Row 33 - was RenderToken, now OKB
Set-TextValue "B33" "OKB"
Set-TextValue "C33" "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
Set-TextValue "D33" "52.97"
Set-TextValue "E33" "  +1.79%  "

# Row 34 - was OKB, now RenderToken
Set-TextValue "B34" "RenderToken"
Set-TextValue "C34" "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-TextValue "D34" "5.94"
Set-TextValue "E34" "  +7.02%  "

# Row 35 - Hedera
Set-TextValue "E35" "  +10.11%  "

# Row 36 - VeChain
Set-TextValue "D36" "0.0460"
Set-TextValue "E36" "  +3.24%  "

# Row 37 - FirstDigitalUSD
Set-TextValue "E37" "  -0.16%  "

# Row 38 - LidoDAOToken
Set-TextValue "D38" "3.29"
Set-TextValue "E38" "  +6.43%  "

# Row 39 - Celestia
Set-TextValue "D39" "18.59"
Set-TextValue "E39" "  -0.27%  "

# Row 40 - ARBITRUM
Set-TextValue "E40" "  +3.62%  "

# Row 41 - Stacks
Set-TextValue "D41" "2.62"
Set-TextValue "E41" "  +4.75%  "

# Row 42 - Stellar
Set-TextValue "E42" "  +2.31%  "

# Row 43 - was Monero, now EnergySwap
Set-TextValue "B43" "EnergySwap"
Set-TextValue "C43" "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue "D43" "22.30"
Set-TextValue "E43" "  +0.88%  "

# Row 44 - was EnergySwap, now Monero
Set-TextValue "B44" "Monero"
Set-TextValue "C44" "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D44" "121.12"
Set-TextValue "E44" "  +0.76%  "

# Row 45 - WEMIXToken
Set-TextValue "D45" "2.21"
Set-TextValue "E45" "  +0.66%  "

# Row 46 - NEARProtocol
Set-TextValue "D46" "3.60"
Set-TextValue "E46" "  +9.65%  "

# Row 47 - Maker
Set-TextValue "D47" "2.176.93"
Set-TextValue "E47" "  +1.66%  "

# Row 48 - ApeXProtocol
Set-TextValue "D48" "2.47"
Set-TextValue "E48" "  +6.31%  "

# Row 49 - TheGraph
Set-TextValue "D49" "0.248"
Set-TextValue "E49" "  +9.62%  "

# Row 50 - SEI
Set-TextValue "D50" "0.961"
Set-TextValue "E50" "  +6.43%  "

# Row 51 - was THORChain, now BEAM
Set-TextValue "B51" "BEAM"
Set-TextValue "C51" "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
Set-TextValue "D51" "0.0323"
Set-TextValue "E51" "  +12.96%  "
